# Update the "6.1.1" indicator metadata sheet with the new responsible-organization
# information (Калымбетова Ы.И. / Управление статистики домашних хозяйств) that
# replaces the previous contact (Керималиева Н.К. / Отдел статистики домашних хозяйств).
#
# The edits are applied in the same order the original author used (B4, then B10,
# B9, B8, B6, B7) so that the workbook's shared-strings table grows in the same
# order, and any now-unused strings are reclaimed, the same way native Excel
# behaves when cell values are overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Indicator title (row 4) - updated wording/punctuation
$ws.Range("B4").Value = "6.1.1 Доля населения, пользующегося услугами водоснабжения, организованного с соблюдением требований безопасности "

# 2. Organization website
$ws.Range("B10").Value = "www.stat.gov.kg"

# 3. Organization phone
$ws.Range("B9").Value = "(0312) 32 46 55"

# 4. Contact person e-mail
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"

# 5. Organization name
$ws.Range("B6").Value = "Национальный статистический комитет КР`n(Управление статистики домашних хозяйств)"

# 6. Contact person name
$ws.Range("B7").Value = "Калымбетова Ы.И."

# Restore the cursor/selection to where the author left it after the last edit
$ws.Activate()
$ws.Range("B10").Select()
